# Updates the "cryptos" price/volume table with refreshed values.
# D/E column literals are prefixed with a leading apostrophe (the
# PowerShell single-quoted string '' ... '' yields a literal leading ')
# so Excel stores them as text (matching the original inlineStr/text
# cells) instead of auto-converting numeric-looking strings to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.239.34'
$ws.Range("E2").Value = '''  +4.89%  '
$ws.Range("D3").Value = '''2.720.21'
$ws.Range("E3").Value = '''  +4.44%  '
$ws.Range("E4").Value = '''  +0.23%  '
$ws.Range("D5").Value = '''586.82'
$ws.Range("E5").Value = '''  +0.49%  '
$ws.Range("D6").Value = '''150.50'
$ws.Range("E6").Value = '''  +5.34%  '
$ws.Range("D7").Value = '''0.997'
$ws.Range("E7").Value = '''  -0.17%  '
$ws.Range("D8").Value = '''0.606'
$ws.Range("E8").Value = '''  +1.79%  '
$ws.Range("D9").Value = '''2.750.97'
$ws.Range("E9").Value = '''  +5.37%  '
$ws.Range("D10").Value = '''6.75'
$ws.Range("E10").Value = '''  +3.68%  '
$ws.Range("E11").Value = '''  +7.28%  '
$ws.Range("D12").Value = '''0.388'
$ws.Range("E12").Value = '''  +4.20%  '
$ws.Range("E13").Value = '''  +1.63%  '
$ws.Range("D14").Value = '''3.236.03'
$ws.Range("E14").Value = '''  +5.60%  '
$ws.Range("D15").Value = '''26.61'
$ws.Range("E15").Value = '''  +8.04%  '
$ws.Range("D16").Value = '''63.190.12'
$ws.Range("E16").Value = '''  +4.85%  '
$ws.Range("D17").Value = '''0.0000150'
$ws.Range("E17").Value = '''  +6.82%  '
$ws.Range("D18").Value = '''2.744.13'
$ws.Range("E18").Value = '''  +5.16%  '
$ws.Range("D19").Value = '''11.97'
$ws.Range("E19").Value = '''  +5.27%  '
$ws.Range("E20").Value = '''  +5.34%  '
$ws.Range("D21").Value = '''364.02'
$ws.Range("E21").Value = '''  +4.81%  '
$ws.Range("D22").Value = '''7.00'
$ws.Range("E22").Value = '''  +1.38%  '
$ws.Range("D23").Value = '''1.00'
$ws.Range("E23").Value = '''  +0.19%  '
$ws.Range("D24").Value = '''0.538'
$ws.Range("E24").Value = '''  +0.19%  '
$ws.Range("D25").Value = '''65.71'
$ws.Range("E25").Value = '''  +2.82%  '
$ws.Range("E26").Value = '''  +4.22%  '
$ws.Range("D27").Value = '''8.61'
$ws.Range("E27").Value = '''  +8.05%  '
$ws.Range("D28").Value = '''0.997'
$ws.Range("E28").Value = '''  +0.01%  '
$ws.Range("D29").Value = '''0.0₃0863'
$ws.Range("E29").Value = '''  +8.43%  '
$ws.Range("E30").Value = '''  +6.90%  '
$ws.Range("D31").Value = '''7.12'
$ws.Range("E31").Value = '''  +11.37%  '
$ws.Range("D32").Value = '''169.80'
$ws.Range("E32").Value = '''  +0.37%  '
$ws.Range("D33").Value = '''0.998'
$ws.Range("E33").Value = '''  -0.05%  '
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D34").Value = '''20.56'
$ws.Range("E34").Value = '''  +5.84%  '
$ws.Range("B35").Value = 'Fetch.AI'
$ws.Range("C35").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D35").Value = '''1.18'
$ws.Range("E35").Value = '''  +18.77%  '
$ws.Range("E36").Value = '''  +11.64%  '
$ws.Range("E37").Value = '''  +9.98%  '
$ws.Range("E38").Value = '''  +9.84%  '
$ws.Range("D39").Value = '''1.02'
$ws.Range("E39").Value = '''  +19.53%  '
$ws.Range("D40").Value = '''348.20'
$ws.Range("E40").Value = '''  +9.98%  '
$ws.Range("E41").Value = '''  +9.72%  '
$ws.Range("D42").Value = '''39.26'
$ws.Range("E42").Value = '''  +2.63%  '
$ws.Range("D43").Value = '''5.66'
$ws.Range("E43").Value = '''  +13.33%  '
$ws.Range("D44").Value = '''22.33'
$ws.Range("E44").Value = '''  +12.10%  '
$ws.Range("D45").Value = '''141.99'
$ws.Range("E45").Value = '''  +4.68%  '
$ws.Range("D46").Value = '''21.90'
$ws.Range("E46").Value = '''  +9.68%  '
$ws.Range("E47").Value = '''  +8.18%  '
$ws.Range("E48").Value = '''  +5.56%  '
$ws.Range("E49").Value = '''  +7.15%  '
$ws.Range("E50").Value = '''  +1.65%  '
$ws.Range("D51").Value = '''2.158.13'
$ws.Range("E51").Value = '''  +6.64%  '
